# Weekly update: a new price observation is inserted as row 289, pushing the
# existing rows 289-362 down to 290-363 (dimension grows from R362 to R363).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 289, shifting everything below it down.
$ws.Rows("289:289").Insert()

# Populate the newly inserted row with the new observation. The columns that
# are constant for every record in this sheet (market/region/category/etc.)
# are copied from the surrounding rows; the observation-specific columns get
# the new values.
$ws.Range("A289").Value = 4
$ws.Range("B289").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C289").Value = "Los Lagos"
$ws.Range("D289").Value = 44754
$ws.Range("E289").Value = 10
$ws.Range("F289").Value = 100112023
$ws.Range("G289").Value = "Brócoli"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 1000
$ws.Range("K289").Value = 1500
$ws.Range("L289").Value = 1500
$ws.Range("M289").Value = 1500
$ws.Range("N289").Value = "$/unidad"
$ws.Range("O289").Value = "Región Metropolitana"
$ws.Range("P289").Value = 1500
$ws.Range("Q289").Value = 1
$ws.Range("R289").Value = "Hortaliza"
